# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
#
# The source feed re-attributed several fixture rows to the wrong match
# record. This script fixes the mis-ordering by re-writing, for each
# affected row, the full set of match columns (B: match id, E: HomeTeam,
# F: AwayTeam, G:AD: score/odds data) using the correct match's data,
# while leaving column A (rank #), C (league), and D (date) untouched
# since those are tied to the row's position, not to a specific match.
#
# Rows 143/144/145 form a 3-way rotation (145 -> 143 -> 144 -> 145).
# Rows 211/212 and 214/215 are simple 2-way swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AD that travel together with a match record (A, C, D stay put)
$cols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Correct (post-fix) data for each affected row, in $cols order.
$rowData = @{
    "143" = @(7532414, "Independiente Petrolero", "Real Santa Cruz", 1, 0, 1, 0, "H", 1.571, 3.75, 5, 1.3, 5, 11, -1.75, 2, 1.8, 3, 1.85, 1.95, 0.3, -1, -1, -1, 0.8, -1, 0.95)
    "144" = @(7532412, "Vaca Diez", "Blooming", 0, 3, 0, 2, "A", 1.727, 3.75, 4, 2.3, 3.6, 2.875, -0.25, 1.95, 1.85, 2.75, 1.925, 1.875, -1, -1, 1.875, -1, 0.8500000000000001, 0.4625, -0.5)
    "145" = @(7532413, "Libertad Gran Mamore FC", "Club Aurora", 0, 1, 0, 0, "A", 2.25, 3.3, 2.8, 2.375, 3.4, 2.875, -0.25, 2.025, 1.775, 2.5, 1.9, 1.9, -1, -1, 1.875, -1, 0.7749999999999999, -1, 0.8999999999999999)
    "211" = @(8039389, "Royal Pari FC", "Club Aurora", 0, 0, 0, 0, "D", 2.2, 3.3, 2.875, 2.625, 3.4, 2.625, 0, 1.875, 1.925, 2.5, 1.95, 1.85, -1, 2.4, -1, 0, 0, -1, 0.8500000000000001)
    "212" = @(8039390, "Blooming", "Nacional Potosi", 2, 1, 0, 1, "H", 2.05, 3.6, 3, 1.8, 3.75, 4.2, -0.5, 1.825, 1.975, 3, 1.975, 1.825, 0.8, -1, -1, 0.825, -1, 0, 0)
    "214" = @(8038943, "San Jose de Oruro", "Bolivar", 2, 1, 1, 0, "H", 2.3, 3.5, 2.625, 2.8, 3.6, 2.375, 0.25, 1.8, 2, 3.25, 1.975, 1.825, 1.8, -1, -1, 0.8, -1, -0.5, 0.4125)
    "215" = @(8039392, "Oriente Petrolero", "Jorge Wilstermann", 2, 1, 0, 0, "H", 2, 3.25, 3.4, 1.727, 4, 4.5, -0.75, 1.9, 1.9, 2.75, 1.9, 1.9, 0.7270000000000001, -1, -1, 0.45, -0.5, 0.45, -0.5)
}

foreach ($rowNum in $rowData.Keys) {
    $values = $rowData[$rowNum]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowNum).Value2 = $values[$i]
    }
}
